$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.08

# Row 3 updates
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48

# Row 5 updates
$ws.Range("G5").Value = 3.75
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 2.1
$ws.Range("J5").Value = 4.5
$ws.Range("L5").Value = 2.88
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5
$ws.Range("AI5").Value = 8.5
$ws.Range("AO5").Value = 23
$ws.Range("AX5").Value = 12
